$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("Upload"), shifting Upload/Desc/
# the data row down by one. Excel auto-extends the data validation ranges
# that already touched row 8 (A6:A8 -> A6:A9, B7:J8 -> B7:J9).
$ws.Rows.Item(8).Insert()

# Give the new row the same look as the row it now sits above (row 7) /
# below (row 9, the old "Upload" row): copy formats from A9:I9 into A8:I8.
$ws.Range("A9:I9").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row content: label "Force" plus all-FALSE flags, matching the layout
# used by the other boolean-flag rows (Save/Cache/Public/Private/Ref/Upload).
$ws.Cells.Item(8, 1).Value2 = "Force"
for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item(8, $c).Value2 = $false
}

# Restore the frozen pane / selection to track the inserted row, same as
# Excel would do once the user re-freezes under the new last header row and
# clicks back on the now-empty row 9.
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("A11").Select() | Out-Null
$aw.FreezePanes = $true
$ws.Range("A9").Select() | Out-Null
